$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Umschaltung AFR <> Lambda"
$ws.Range("A2").Value = "Cursor-Werte"
$ws.Range("A3").Value = "Print graph"
$ws.Range("A4").Value = "Automatenmodus"

$ws.Range("A5").Select()
